# Update "想去人数" (want-to-go count) and one "最低票价" (min price) figure
# across the two sheets that carry this data: "展览" and "全部类型".
# Values mirror each other row-for-row between the two sheets (same events),
# but the row numbers differ because "全部类型" has one extra row inserted
# above the affected rows later in the list.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 6493
$wsExhibit.Range("F4").Value = 8
$wsExhibit.Range("F5").Value = 397
$wsExhibit.Range("F10").Value = 81
$wsExhibit.Range("F12").Value = 160
$wsExhibit.Range("G13").Value = 55
$wsExhibit.Range("F15").Value = 3191
$wsExhibit.Range("F16").Value = 15
$wsExhibit.Range("F18").Value = 1863

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 6493
$wsAll.Range("F4").Value = 8
$wsAll.Range("F5").Value = 397
$wsAll.Range("F11").Value = 81
$wsAll.Range("F13").Value = 160
$wsAll.Range("G14").Value = 55
$wsAll.Range("F16").Value = 3191
$wsAll.Range("F17").Value = 15
$wsAll.Range("F19").Value = 1863
